# Apply text ("@") number format to Price cells whose new values look numeric,
# so Excel stores them as text (matching the source data which uses inline strings)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5:D6").NumberFormat = "@"
$ws.Range("D8:D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19:D36").NumberFormat = "@"
$ws.Range("D38:D39").NumberFormat = "@"
$ws.Range("D42:D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.858.59"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "1.633.22"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("D5").Value = "213.75"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").Value = "0.5057"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").Value = "0.2567"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "0.06342"
$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").Value = "0.07729"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").Value = "4.261"
$ws.Range("E12").Value = "  -0.71%  "

$ws.Range("D13").Value = "1.633.02"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").Value = "0.5429"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "0.0₅7708"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "63.99"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "25.883.29"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "195.10"
$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "4.418"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").Value = "9.894"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").Value = "6.090"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").Value = "1.891"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").Value = "142.72"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").Value = "0.1246"
$ws.Range("E26").Value = "  +7.29%  "

$ws.Range("D27").Value = "6.796"
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("D28").Value = "15.59"
$ws.Range("E28").Value = "  -1.08%  "

$ws.Range("D29").Value = "1.234"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").Value = "0.04848"
$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").Value = "3.230"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").Value = "3.177"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").Value = "1.541"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").Value = "2.372"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").Value = "0.9055"
$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("D36").Value = "2.573"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").Value = "1.125.73"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").Value = "0.5478"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").Value = "0.01553"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").Value = "0.8003"
$ws.Range("E42").Value = "  -2.00%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "98.34"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₈123"
$ws.Range("E44").Value = "  -5.94%  "

$ws.Range("D45").Value = "1.770.87"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "0.4474"
$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").Value = "54.83"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "0.05161"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").Value = "7.535"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.55%  "
